$d = $word.ActiveDocument

function Set-ParagraphText($index, $text) {
    $p = $d.Paragraphs.Item($index)
    $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
    $rng.Text = $text
}

# Paragraph 2: Studente(...) - move Email after Username
Set-ParagraphText 2 "Studente(Matricola(PK), Username, Email, Nome, Cognome, Passkey)"

# Paragraph 3: Oggetto(...) - new field list
Set-ParagraphText 3 "Oggetto(idOggetto(PK), Descrizione, NomeCategoria(FK), Matricola(FK))"

# Paragraph 4 (was Annuncio): becomes the new Categoria table
Set-ParagraphText 4 "Categoria(NomeCategoria(PK))"

# Paragraph 5 (was Sede): becomes the new Annuncio table
Set-ParagraphText 5 "Annuncio(idAnnuncio(PK), Titolo, StatoAnnuncio, FasciaOrariaInizio, FasciaOrariaFine, Descrizione, Prezzo, Tipologia, DataPubblicazione, Matricola(FK), idOggetto(FK), idSede(FK))"

# Paragraph 6 (was Offerta): becomes the new Sede table
Set-ParagraphText 6 "Sede(idSede(PK), Ptop, Descrizione, Civico, CAP)"

# Paragraph 7 (was OggettoOfferto): becomes the new Offerta table
Set-ParagraphText 7 "Offerta(idOfferta(PK), StatoOfferta, PrezzoOfferta, Motivazione, Tipologia, DataInvio, Matricola(FK), idAnnuncio(FK))"

# Append a brand-new final paragraph for OggettoOfferto
$lastPara = $d.Paragraphs.Item(7)
$lastPara.Range.InsertParagraphAfter()
Set-ParagraphText 8 "OggettoOfferto(IdOfferta(FK), IdOggetto(FK))"
